$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells hold plain text in the source data (coin names, links,
# formatted price strings, and padded percentage strings). Force the
# number format to Text first so Excel does not auto-convert numeric-
# looking strings (e.g. "89.69", "1.80") into floating point numbers,
# which would silently change/round the displayed text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.167.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.432.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.70%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0838"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.15"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.803.81"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.400.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.098.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0929"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.56%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.63"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.68"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.73"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0748"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.97"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.96"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.31"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.25%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.996.29"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.04%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.666.61"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.60"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.90"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.64%  "
